$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.991.01'
$ws.Range('E2').Value = '  +0.51%  '
$ws.Range('D3').Value = '1.641.11'
$ws.Range('E3').Value = '  +0.52%  '
$ws.Range('E4').Value = '  +0.42%  '
$ws.Range('D5').Value = '215.83'
$ws.Range('E5').Value = '  +0.66%  '
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('E7').Value = '  +0.41%  '
$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D8').Value = '0.255'
$ws.Range('E8').Value = '  +0.58%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').Value = '0.0638'
$ws.Range('E9').Value = '  +0.88%  '
$ws.Range('D10').Value = '19.52'
$ws.Range('E10').Value = '  -0.13%  '
$ws.Range('D11').Value = '0.0796'
$ws.Range('E11').Value = '  +0.58%  '
$ws.Range('D12').Value = '1.869.42'
$ws.Range('E12').Value = '  +0.64%  '
$ws.Range('E13').Value = '  +0.60%  '
$ws.Range('D14').Value = '1.656.76'
$ws.Range('E14').Value = '  +0.89%  '
$ws.Range('E15').Value = '  -0.07%  '
$ws.Range('D16').Value = '0.0₃0763'
$ws.Range('E16').Value = '  +0.90%  '
$ws.Range('D17').Value = '63.37'
$ws.Range('E17').Value = '  +1.33%  '
$ws.Range('D18').Value = '26.099.07'
$ws.Range('E18').Value = '  +0.93%  '
$ws.Range('E19').Value = '  +0.38%  '
$ws.Range('D20').Value = '194.06'
$ws.Range('E20').Value = '  +0.29%  '
$ws.Range('E21').Value = '  -0.68%  '
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('E23').Value = '  -1.20%  '
$ws.Range('B24').Value = 'Stellar'
$ws.Range('C24').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D24').Value = '0.131'
$ws.Range('E24').Value = '  +4.80%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Value = '1.80'
$ws.Range('E25').Value = '  -1.11%  '
$ws.Range('E26').Value = '  +0.37%  '
$ws.Range('D27').Value = '142.97'
$ws.Range('E27').Value = '  -0.32%  '
$ws.Range('E28').Value = '  +0.50%  '
$ws.Range('E29').Value = '  +0.71%  '
$ws.Range('E30').Value = '  +0.70%  '
$ws.Range('E31').Value = '  -0.66%  '
$ws.Range('E32').Value = '  -0.35%  '
$ws.Range('E33').Value = '  +1.12%  '
$ws.Range('E34').Value = '  -0.97%  '
$ws.Range('E35').Value = '  +1.29%  '
$ws.Range('E36').Value = '  +0.35%  '
$ws.Range('D37').Value = '1.130.40'
$ws.Range('E37').Value = '  -0.70%  '
$ws.Range('D38').Value = '0.539'
$ws.Range('E38').Value = '  -1.08%  '
$ws.Range('E39').Value = '  -0.49%  '
$ws.Range('D40').Value = '0.0157'
$ws.Range('E40').Value = '  +0.24%  '
$ws.Range('D42').Value = '99.12'
$ws.Range('E42').Value = '  -0.15%  '
$ws.Range('E43').Value = '  -0.11%  '
$ws.Range('D44').Value = '1.778.60'
$ws.Range('E44').Value = '  +0.68%  '
$ws.Range('E45').Value = '  +4.95%  '
$ws.Range('E46').Value = '  +0.49%  '
$ws.Range('E48').Value = '  +3.18%  '
$ws.Range('D49').Value = '7.72'
$ws.Range('E49').Value = '  +1.14%  '
$ws.Range('E50').Value = '  -0.15%  '
$ws.Range('E51').Value = '  +0.27%  '
